$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -12.998
$ws.Range("B7").Value = 5.605
$ws.Range("E7").Value = 16.069
$ws.Range("A9").Value = -21.864
$ws.Range("E10").Value = 16.437
$ws.Range("B12").Value = 5.662
$ws.Range("A13").Value = -21.979
$ws.Range("E13").Value = 16.477
$ws.Range("B14").Value = 6.044
$ws.Range("C15").Value = -13.134
$ws.Range("A16").Value = -21.83
$ws.Range("E16").Value = 16.625
$ws.Range("A18").Value = -21.985
$ws.Range("B19").Value = 7.889
$ws.Range("A20").Value = -20.846
$ws.Range("E20").Value = 16.452
$ws.Range("E24").Value = 17.006
$ws.Range("A26").Value = -21.269
$ws.Range("B26").Value = 6.734999999999999
$ws.Range("A27").Value = -21.711
$ws.Range("B27").Value = 5.829
$ws.Range("C28").Value = -12.965
$ws.Range("A29").Value = -21.543
$ws.Range("B29").Value = 6.437
$ws.Range("E32").Value = 16.15
$ws.Range("C33").Value = -11.292
$ws.Range("A35").Value = -20.452
$ws.Range("C35").Value = -12.624
$ws.Range("A36").Value = -21.089
$ws.Range("B37").Value = 7.515000000000001
$ws.Range("B38").Value = 5.557
$ws.Range("C38").Value = -12.647
$ws.Range("E39").Value = 16.222
$ws.Range("C43").Value = -12.774
$ws.Range("C44").Value = -12.277
$ws.Range("A45").Value = -21.451
$ws.Range("C45").Value = -13.244
$ws.Range("B47").Value = 5.635999999999999
$ws.Range("C47").Value = -13.038
$ws.Range("E47").Value = 16.585
$ws.Range("E48").Value = 17.19
$ws.Range("B51").Value = 5.816
$ws.Range("C51").Value = -11.76
$ws.Range("B52").Value = 5.678
$ws.Range("E52").Value = 16.925
$ws.Range("C54").Value = -13.011
$ws.Range("A55").Value = -21.632
$ws.Range("B55").Value = 6.392
$ws.Range("E56").Value = 16.656
$ws.Range("A57").Value = -21.337
$ws.Range("C57").Value = -13.051
$ws.Range("C62").Value = -13.409
$ws.Range("C63").Value = -12.484
$ws.Range("C67").Value = -11.224
$ws.Range("A69").Value = -21.538
$ws.Range("B69").Value = 6.256
$ws.Range("B70").Value = 6.08
$ws.Range("C70").Value = -11.114
$ws.Range("A76").Value = -21.706
$ws.Range("B76").Value = 6.285
$ws.Range("A78").Value = -20.427
$ws.Range("B81").Value = 4.685
$ws.Range("C81").Value = -13.317
$ws.Range("A82").Value = -21.875
$ws.Range("A83").Value = -20.659
$ws.Range("B83").Value = 7.124000000000001
$ws.Range("E84").Value = 16.663
$ws.Range("C88").Value = -13.25
$ws.Range("A93").Value = -21.573
$ws.Range("B94").Value = 6.616999999999999
$ws.Range("C96").Value = -12.998
$ws.Range("A97").Value = -21.769
$ws.Range("C99").Value = -12.984
$ws.Range("B100").Value = 6.124
$ws.Range("E100").Value = 16.476
$ws.Range("E101").Value = 16.793
$ws.Range("B102").Value = 6.948
